$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col14a1"
$ws.Cells.Item(2, 3).Value = "Cd44"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.5688816666666666
$ws.Cells.Item(2, 8).Value = 1.706645
$ws.Cells.Item(2, 9).Value = 0.001790814942693061
$ws.Cells.Item(2, 10).Value = 0.001790814942693061
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 31.82741333333333
$ws.Cells.Item(2, 14).Value = 95.48223999999999
$ws.Cells.Item(2, 15).Value = 0.114390792932228
$ws.Cells.Item(2, 16).Value = 0.114390792932228
$ws.Cells.Item(2, 17).Value = 18.10603194275555
$ws.Cells.Item(2, 18).Value = 162.9542874848
$ws.Cells.Item(2, 19).Value = 0.0002048527412895417
$ws.Cells.Item(2, 20).Value = 0.0002048527412895418

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col14a1"
$ws.Cells.Item(3, 3).Value = "Cd44"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.5688816666666666
$ws.Cells.Item(3, 8).Value = 1.706645
$ws.Cells.Item(3, 9).Value = 0.001790814942693061
$ws.Cells.Item(3, 10).Value = 0.001790814942693061
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 85.46317833333335
$ws.Cells.Item(3, 14).Value = 256.389535
$ws.Cells.Item(3, 15).Value = 0.307162904935779
$ws.Cells.Item(3, 16).Value = 0.307162904935779
$ws.Cells.Item(3, 17).Value = 48.61843532889723
$ws.Cells.Item(3, 18).Value = 437.565917960075
$ws.Cells.Item(3, 19).Value = 0.0005500719200000014
$ws.Cells.Item(3, 20).Value = 0.0005500719200000014

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col14a1"
$ws.Cells.Item(4, 3).Value = "Cd44"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.5688816666666666
$ws.Cells.Item(4, 8).Value = 1.706645
$ws.Cells.Item(4, 9).Value = 0.001790814942693061
$ws.Cells.Item(4, 10).Value = 0.001790814942693061
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 122.2478306666667
$ws.Cells.Item(4, 14).Value = 366.743492
$ws.Cells.Item(4, 15).Value = 0.4393704929064738
$ws.Cells.Item(4, 16).Value = 0.4393704929064738
$ws.Cells.Item(4, 17).Value = 69.54454965603777
$ws.Cells.Item(4, 18).Value = 625.90094690434
$ws.Cells.Item(4, 19).Value = 0.000786831244075329
$ws.Cells.Item(4, 20).Value = 0.0007868312440753291

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Col14a1"
$ws.Cells.Item(5, 3).Value = "Cd44"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.5688816666666666
$ws.Cells.Item(5, 8).Value = 1.706645
$ws.Cells.Item(5, 9).Value = 0.001790814942693061
$ws.Cells.Item(5, 10).Value = 0.001790814942693061
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 38.69562533333333
$ws.Cells.Item(5, 14).Value = 116.086876
$ws.Cells.Item(5, 15).Value = 0.1390758092255191
$ws.Cells.Item(5, 16).Value = 0.1390758092255191
$ws.Cells.Item(5, 17).Value = 22.01323183233555
$ws.Cells.Item(5, 18).Value = 198.11908649102
$ws.Cells.Item(5, 19).Value = 0.0002490590373281892
$ws.Cells.Item(5, 20).Value = 0.0002490590373281892

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col14a1"
$ws.Cells.Item(6, 3).Value = "Cd44"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 314.9820043333334
$ws.Cells.Item(6, 8).Value = 944.946013
$ws.Cells.Item(6, 9).Value = 0.9915497599785732
$ws.Cells.Item(6, 10).Value = 0.991549759978573
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 31.82741333333333
$ws.Cells.Item(6, 14).Value = 95.48223999999999
$ws.Cells.Item(6, 15).Value = 0.114390792932228
$ws.Cells.Item(6, 16).Value = 0.114390792932228
$ws.Cells.Item(6, 17).Value = 10025.06244447879
$ws.Cells.Item(6, 18).Value = 90225.56200030912
$ws.Cells.Item(6, 19).Value = 0.1134241632757094
$ws.Cells.Item(6, 20).Value = 0.1134241632757094

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col14a1"
$ws.Cells.Item(7, 3).Value = "Cd44"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 314.9820043333334
$ws.Cells.Item(7, 8).Value = 944.946013
$ws.Cells.Item(7, 9).Value = 0.9915497599785732
$ws.Cells.Item(7, 10).Value = 0.991549759978573
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 85.46317833333335
$ws.Cells.Item(7, 14).Value = 256.389535
$ws.Cells.Item(7, 15).Value = 0.307162904935779
$ws.Cells.Item(7, 16).Value = 0.307162904935779
$ws.Cells.Item(7, 17).Value = 26919.36320813045
$ws.Cells.Item(7, 18).Value = 242274.268873174
$ws.Cells.Item(7, 19).Value = 0.304567304663393
$ws.Cells.Item(7, 20).Value = 0.3045673046633929

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Col14a1"
$ws.Cells.Item(8, 3).Value = "Cd44"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 314.9820043333334
$ws.Cells.Item(8, 8).Value = 944.946013
$ws.Cells.Item(8, 9).Value = 0.9915497599785732
$ws.Cells.Item(8, 10).Value = 0.991549759978573
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 122.2478306666667
$ws.Cells.Item(8, 14).Value = 366.743492
$ws.Cells.Item(8, 15).Value = 0.4393704929064738
$ws.Cells.Item(8, 16).Value = 0.4393704929064738
$ws.Cells.Item(8, 17).Value = 38505.86672878861
$ws.Cells.Item(8, 18).Value = 346552.8005590974
$ws.Cells.Item(8, 19).Value = 0.4356577067830815
$ws.Cells.Item(8, 20).Value = 0.4356577067830814

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Col14a1"
$ws.Cells.Item(9, 3).Value = "Cd44"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 314.9820043333334
$ws.Cells.Item(9, 8).Value = 944.946013
$ws.Cells.Item(9, 9).Value = 0.9915497599785732
$ws.Cells.Item(9, 10).Value = 0.991549759978573
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 38.69562533333333
$ws.Cells.Item(9, 14).Value = 116.086876
$ws.Cells.Item(9, 15).Value = 0.1390758092255191
$ws.Cells.Item(9, 16).Value = 0.1390758092255191
$ws.Cells.Item(9, 17).Value = 12188.42562642504
$ws.Cells.Item(9, 18).Value = 109695.8306378254
$ws.Cells.Item(9, 19).Value = 0.1379005852563893
$ws.Cells.Item(9, 20).Value = 0.1379005852563893

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Col14a1"
$ws.Cells.Item(10, 3).Value = "Cd44"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.73284
$ws.Cells.Item(10, 8).Value = 5.19852
$ws.Cells.Item(10, 9).Value = 0.005454905557915521
$ws.Cells.Item(10, 10).Value = 0.005454905557915521
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 31.82741333333333
$ws.Cells.Item(10, 14).Value = 95.48223999999999
$ws.Cells.Item(10, 15).Value = 0.114390792932228
$ws.Cells.Item(10, 16).Value = 0.114390792932228
$ws.Cells.Item(10, 17).Value = 55.15181492053333
$ws.Cells.Item(10, 18).Value = 496.3663342848
$ws.Cells.Item(10, 19).Value = 0.0006239909721403741
$ws.Cells.Item(10, 20).Value = 0.0006239909721403741

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Col14a1"
$ws.Cells.Item(11, 3).Value = "Cd44"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.73284
$ws.Cells.Item(11, 8).Value = 5.19852
$ws.Cells.Item(11, 9).Value = 0.005454905557915521
$ws.Cells.Item(11, 10).Value = 0.005454905557915521
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 85.46317833333335
$ws.Cells.Item(11, 14).Value = 256.389535
$ws.Cells.Item(11, 15).Value = 0.307162904935779
$ws.Cells.Item(11, 16).Value = 0.307162904935779
$ws.Cells.Item(11, 17).Value = 148.0940139431334
$ws.Cells.Item(11, 18).Value = 1332.8461254882
$ws.Cells.Item(11, 19).Value = 0.001675544637319658
$ws.Cells.Item(11, 20).Value = 0.001675544637319658

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Col14a1"
$ws.Cells.Item(12, 3).Value = "Cd44"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.73284
$ws.Cells.Item(12, 8).Value = 5.19852
$ws.Cells.Item(12, 9).Value = 0.005454905557915521
$ws.Cells.Item(12, 10).Value = 0.005454905557915521
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 122.2478306666667
$ws.Cells.Item(12, 14).Value = 366.743492
$ws.Cells.Item(12, 15).Value = 0.4393704929064738
$ws.Cells.Item(12, 16).Value = 0.4393704929064738
$ws.Cells.Item(12, 17).Value = 211.8359308924267
$ws.Cells.Item(12, 18).Value = 1906.52337803184
$ws.Cells.Item(12, 19).Value = 0.002396724543739606
$ws.Cells.Item(12, 20).Value = 0.002396724543739606

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Col14a1"
$ws.Cells.Item(13, 3).Value = "Cd44"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.73284
$ws.Cells.Item(13, 8).Value = 5.19852
$ws.Cells.Item(13, 9).Value = 0.005454905557915521
$ws.Cells.Item(13, 10).Value = 0.005454905557915521
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 38.69562533333333
$ws.Cells.Item(13, 14).Value = 116.086876
$ws.Cells.Item(13, 15).Value = 0.1390758092255191
$ws.Cells.Item(13, 16).Value = 0.1390758092255191
$ws.Cells.Item(13, 17).Value = 67.05332740261333
$ws.Cells.Item(13, 18).Value = 603.4799466235199
$ws.Cells.Item(13, 19).Value = 0.000758645404715883
$ws.Cells.Item(13, 20).Value = 0.0007586454047158829

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Col14a1"
$ws.Cells.Item(14, 3).Value = "Cd44"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.3826353333333333
$ws.Cells.Item(14, 8).Value = 1.147906
$ws.Cells.Item(14, 9).Value = 0.001204519520818343
$ws.Cells.Item(14, 10).Value = 0.001204519520818343
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 31.82741333333333
$ws.Cells.Item(14, 14).Value = 95.48223999999999
$ws.Cells.Item(14, 15).Value = 0.114390792932228
$ws.Cells.Item(14, 16).Value = 0.114390792932228
$ws.Cells.Item(14, 17).Value = 12.17829290993777
$ws.Cells.Item(14, 18).Value = 109.60463618944
$ws.Cells.Item(14, 19).Value = 0.0001377859430887576
$ws.Cells.Item(14, 20).Value = 0.0001377859430887576

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Col14a1"
$ws.Cells.Item(15, 3).Value = "Cd44"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.3826353333333333
$ws.Cells.Item(15, 8).Value = 1.147906
$ws.Cells.Item(15, 9).Value = 0.001204519520818343
$ws.Cells.Item(15, 10).Value = 0.001204519520818343
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 85.46317833333335
$ws.Cells.Item(15, 14).Value = 256.389535
$ws.Cells.Item(15, 15).Value = 0.307162904935779
$ws.Cells.Item(15, 16).Value = 0.307162904935779
$ws.Cells.Item(15, 17).Value = 32.70123172930111
$ws.Cells.Item(15, 18).Value = 294.31108556371
$ws.Cells.Item(15, 19).Value = 0.0003699837150664148
$ws.Cells.Item(15, 20).Value = 0.0003699837150664148

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Col14a1"
$ws.Cells.Item(16, 3).Value = "Cd44"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.3826353333333333
$ws.Cells.Item(16, 8).Value = 1.147906
$ws.Cells.Item(16, 9).Value = 0.001204519520818343
$ws.Cells.Item(16, 10).Value = 0.001204519520818343
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 122.2478306666667
$ws.Cells.Item(16, 14).Value = 366.743492
$ws.Cells.Item(16, 15).Value = 0.4393704929064738
$ws.Cells.Item(16, 16).Value = 0.4393704929064738
$ws.Cells.Item(16, 17).Value = 46.77633943641688
$ws.Cells.Item(16, 18).Value = 420.987054927752
$ws.Cells.Item(16, 19).Value = 0.000529230335577425
$ws.Cells.Item(16, 20).Value = 0.000529230335577425

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Col14a1"
$ws.Cells.Item(17, 3).Value = "Cd44"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.3826353333333333
$ws.Cells.Item(17, 8).Value = 1.147906
$ws.Cells.Item(17, 9).Value = 0.001204519520818343
$ws.Cells.Item(17, 10).Value = 0.001204519520818343
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 38.69562533333333
$ws.Cells.Item(17, 14).Value = 116.086876
$ws.Cells.Item(17, 15).Value = 0.1390758092255191
$ws.Cells.Item(17, 16).Value = 0.1390758092255191
$ws.Cells.Item(17, 17).Value = 14.80631349796177
$ws.Cells.Item(17, 18).Value = 133.256821481656
$ws.Cells.Item(17, 19).Value = 0.0001675195270857456
$ws.Cells.Item(17, 20).Value = 0.0001675195270857456

